$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.540.10"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "3.445.62"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'577.75"

$ws.Range("D6").Value = "'145.32"
$ws.Range("E6").Value = "  +4.40%  "

$ws.Range("D7").Value = "3.446.32"
$ws.Range("E7").Value = "  +1.17%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.480"
$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("E11").Value = "  +3.54%  "

$ws.Range("D12").Value = "'0.390"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").Value = "4.033.40"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").Value = "'28.46"
$ws.Range("E14").Value = "  +6.46%  "

$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").Value = "3.447.10"
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").Value = "61.648.15"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "  +6.73%  "

$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = "  +3.51%  "

$ws.Range("D21").Value = "'9.46"
$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("D22").Value = "'403.41"
$ws.Range("E22").Value = "  +6.59%  "

$ws.Range("D23").Value = "'0.570"
$ws.Range("E23").Value = "  +3.12%  "

$ws.Range("D24").Value = "'74.46"
$ws.Range("E24").Value = "  +4.45%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.77"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("D27").Value = "'0.0000124"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").Value = "3.587.10"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("E29").Value = "  +4.46%  "

$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  -10.43%  "

$ws.Range("D36").Value = "'23.92"
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("D37").Value = "'7.06"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").Value = "3.473.50"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").Value = "'1.57"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").Value = "'167.09"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("D43").Value = "'27.25"
$ws.Range("E43").Value = "  +3.72%  "

$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("E45").Value = "  +3.29%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.74"
$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D49").Value = "2.613.54"
$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("E51").Value = "  +2.43%  "
